$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# New identifiers / values used throughout the report refresh
# ---------------------------------------------------------------
$newMdA      = "5af33faa-e0ed-4aa5-9e7f-6600e42497e3.md"
$newMdB      = "ffffde93bf6f-127a-4e4a-af53-24c1ce81b6c7.md"
$statusText  = "Ready for handoff"
$overviewDate = "2016-51-18 14:51:32"

$zhXlf        = "5af33faa-e0ed-4aa5-9e7f-6600e42497e3.6e6cabe52da61819d3a3d382cc5527cc49a3f56b.zh-cn.xlf"
$deXlf        = "5af33faa-e0ed-4aa5-9e7f-6600e42497e3.6e6cabe52da61819d3a3d382cc5527cc49a3f56b.de-de.xlf"
$zhHandoffDt  = "2016-03-18 14:51:29"
$deHandoffDt  = "2016-03-18 14:51:32"
$handbackDt   = "0001-01-01 00:00:00"

# =================================================================
# Sheet "Overview"
# =================================================================
$wsO = $wb.Worksheets.Item("Overview")

$wsO.Range("A2").Value = $newMdA
$wsO.Range("B2").Value = $statusText
$wsO.Range("C2").Value = $statusText
$wsO.Range("D2").Value = $overviewDate

$wsO.Range("A3").Value = $newMdB
$wsO.Range("B3").Value = $statusText
$wsO.Range("C3").Value = $statusText
$wsO.Range("D3").Value = $overviewDate

$wsO.Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.md", [Type]::Missing, [Type]::Missing, $newMdA) | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/ffffde93bf6f-127a-4e4a-af53-24c1ce81b6c7.md", [Type]::Missing, [Type]::Missing, $newMdB) | Out-Null

# =================================================================
# Sheet "zh-cn"
# =================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdA
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("D2").Value = $zhXlf
$wsZh.Range("E2").Value = $zhHandoffDt
$wsZh.Range("F2:G2").Clear()
$wsZh.Range("H2").Value = $handbackDt
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = $newMdB
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("D3").Value = $zhXlf
$wsZh.Range("E3").Value = $zhHandoffDt
$wsZh.Range("F3:G3").Clear()
$wsZh.Range("H3").Value = $handbackDt
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.md", [Type]::Missing, [Type]::Missing, $newMdA) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/76d3ffcce4303fe9a7ad70f2626bc3f71b2d8b77/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.6e6cabe52da61819d3a3d382cc5527cc49a3f56b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, $zhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/ffffde93bf6f-127a-4e4a-af53-24c1ce81b6c7.md", [Type]::Missing, [Type]::Missing, $newMdB) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/ffffde93bf6f-127a-4e4a-af53-24c1ce81b6c7.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/76d3ffcce4303fe9a7ad70f2626bc3f71b2d8b77/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.6e6cabe52da61819d3a3d382cc5527cc49a3f56b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, $zhXlf) | Out-Null

# =================================================================
# Sheet "de-de"
# =================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdA
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("D2").Value = $deXlf
$wsDe.Range("E2").Value = $deHandoffDt
$wsDe.Range("F2:G2").Clear()
$wsDe.Range("H2").Value = $handbackDt
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = $newMdB
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("D3").Value = $deXlf
$wsDe.Range("E3").Value = $deHandoffDt
$wsDe.Range("F3:G3").Clear()
$wsDe.Range("H3").Value = $handbackDt
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.md", [Type]::Missing, [Type]::Missing, $newMdA) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25e6e56e3c4a15512659fe5cfb852c4b14313391/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.6e6cabe52da61819d3a3d382cc5527cc49a3f56b.de-de.xlf", [Type]::Missing, [Type]::Missing, $deXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/ffffde93bf6f-127a-4e4a-af53-24c1ce81b6c7.md", [Type]::Missing, [Type]::Missing, $newMdB) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ef917f1cef96bc95c1854f6b20b67bf7da79b4b/e2e/ffffde93bf6f-127a-4e4a-af53-24c1ce81b6c7.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25e6e56e3c4a15512659fe5cfb852c4b14313391/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5af33faa-e0ed-4aa5-9e7f-6600e42497e3.6e6cabe52da61819d3a3d382cc5527cc49a3f56b.de-de.xlf", [Type]::Missing, [Type]::Missing, $deXlf) | Out-Null

$wb.Save()
